$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" "298.07"
Set-TextValue "E2" "-0.21%"
Set-TextValue "D3" "31.27"
Set-TextValue "E3" "0.06%"
Set-TextValue "D4" "5.098"
Set-TextValue "E4" "-0.54%"
Set-TextValue "D5" "0.08037"
Set-TextValue "E5" "9.60%"
Set-TextValue "D6" "2.437"
Set-TextValue "E6" "32.52%"
Set-TextValue "D7" "7.800"
Set-TextValue "E7" "0.67%"
Set-TextValue "D8" "3.797"
Set-TextValue "E8" "1.94%"
Set-TextValue "D9" "0.9183"
Set-TextValue "D10" "0.1727"
Set-TextValue "E10" "3.32%"
Set-TextValue "D11" "0.07337"
Set-TextValue "E11" "3.90%"
Set-TextValue "D12" "0.08575"
Set-TextValue "E12" "7.26%"
Set-TextValue "D13" "0.03026"
Set-TextValue "E13" "1.00%"
Set-TextValue "D14" "0.09968"
Set-TextValue "E14" "0.56%"
Set-TextValue "D15" "0.001508"
Set-TextValue "E15" "1.14%"
Set-TextValue "D16" "0.005976"
Set-TextValue "E16" "-2.72%"
Set-TextValue "D17" "3.510"
Set-TextValue "E17" "1.58%"
Set-TextValue "E18" "0.91%"
Set-TextValue "E19" "1.78%"
Set-TextValue "E20" "1.64%"
Set-TextValue "D21" "4.607"
Set-TextValue "E21" "1.33%"
Set-TextValue "E22" "2.31%"
Set-TextValue "D23" "0.04620"
Set-TextValue "E23" "-0.46%"
Set-TextValue "D24" "0.001247"
Set-TextValue "E24" "2.90%"
Set-TextValue "D25" "0.004438"
Set-TextValue "E25" "-6.38%"
Set-TextValue "E26" "-7.34%"
Set-TextValue "D27" "0.0003427"
Set-TextValue "E27" "83.13%"
Set-TextValue "D39" "0.01798"
Set-TextValue "E39" "4.21%"
Set-TextValue "D40" "0.04512"
Set-TextValue "E40" "0.92%"
Set-TextValue "D41" "0.007073"
Set-TextValue "E41" "-0.73%"
Set-TextValue "D42" "0.1342"
Set-TextValue "E42" "0.87%"
Set-TextValue "D43" "0.002242"
Set-TextValue "E43" "2.21%"
Set-TextValue "D44" "0.009835"
Set-TextValue "E44" "-8.98%"
Set-TextValue "D45" "0.00006601"
Set-TextValue "E45" "5.87%"
Set-TextValue "E46" "-0.02%"
Set-TextValue "E47" "-55.58%"
Set-TextValue "D49" "0.00002100"
Set-TextValue "E49" "-0.02%"
Set-TextValue "D50" "0.0002000"
Set-TextValue "E50" "0.05%"
